$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update column F values for rows identified below
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 293
$ws1.Range("F3").Value = 60
$ws1.Range("F4").Value = 3591
$ws1.Range("F5").Value = 2217
$ws1.Range("F8").Value = 173
$ws1.Range("F9").Value = 80
$ws1.Range("F11").Value = 1326
$ws1.Range("F13").Value = 1910
$ws1.Range("F14").Value = 138

# Sheet "全部类型" (sheet4) - update column F values for the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 293
$ws4.Range("F3").Value = 60
$ws4.Range("F4").Value = 3591
$ws4.Range("F5").Value = 2217
$ws4.Range("F9").Value = 173
$ws4.Range("F10").Value = 80
$ws4.Range("F14").Value = 1326
$ws4.Range("F16").Value = 1910
$ws4.Range("F17").Value = 138
